$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

$ws.Range("D2").Value = "29.602.86"
$ws.Range("E2").Value = "  -0.36%  "
$ws.Range("D3").Value = "1.854.87"
$ws.Range("E3").Value = "  -0.07%  "
Set-TextValue $ws.Range("D4") "0.9992"
$ws.Range("E4").Value = "  -0.14%  "
Set-TextValue $ws.Range("D5") "243.73"
$ws.Range("E5").Value = "  -0.28%  "
Set-TextValue $ws.Range("D6") "0.6407"
$ws.Range("E6").Value = "  +0.43%  "
Set-TextValue $ws.Range("D7") "1.000"
$ws.Range("E7").Value = "  -0.10%  "
Set-TextValue $ws.Range("D8") "0.07591"
$ws.Range("E8").Value = "  +1.51%  "
Set-TextValue $ws.Range("D9") "0.3004"
$ws.Range("E9").Value = "  +0.51%  "
Set-TextValue $ws.Range("D10") "24.42"
$ws.Range("E10").Value = "  +0.74%  "
Set-TextValue $ws.Range("D11") "0.07682"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").Value = "1.867.97"
$ws.Range("E12").Value = "  +0.53%  "
Set-TextValue $ws.Range("D14") "0.6905"
$ws.Range("E14").Value = "  +0.52%  "
Set-TextValue $ws.Range("D15") "84.25"
$ws.Range("E15").Value = "  +0.39%  "
Set-TextValue $ws.Range("D16") "0.000009686"
$ws.Range("E16").Value = "  +1.59%  "
$ws.Range("D17").Value = "2.135.36"
$ws.Range("E17").Value = "  +0.76%  "
Set-TextValue $ws.Range("D18") "6.316"
$ws.Range("E18").Value = "  +4.49%  "
$ws.Range("D19").Value = "29.648.24"
$ws.Range("E19").Value = "  -0.27%  "
Set-TextValue $ws.Range("D20") "239.36"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("E21").Value = "  -0.09%  "
Set-TextValue $ws.Range("D22") "0.9999"
$ws.Range("E22").Value = "  -0.12%  "
Set-TextValue $ws.Range("D23") "7.647"
$ws.Range("E23").Value = "  +3.21%  "
Set-TextValue $ws.Range("D24") "1.000"
$ws.Range("E24").Value = "  -0.13%  "
Set-TextValue $ws.Range("D25") "157.11"
$ws.Range("E25").Value = "  -0.63%  "
Set-TextValue $ws.Range("D26") "0.1405"
$ws.Range("E26").Value = "  -0.93%  "
Set-TextValue $ws.Range("D27") "8.518"
$ws.Range("E27").Value = "  +0.29%  "
Set-TextValue $ws.Range("D28") "17.80"
$ws.Range("E28").Value = "  -0.66%  "
Set-TextValue $ws.Range("D29") "1.492"
$ws.Range("E29").Value = "  +0.11%  "
Set-TextValue $ws.Range("D30") "0.05888"
$ws.Range("E30").Value = "  -4.36%  "
Set-TextValue $ws.Range("D31") "1.286"
$ws.Range("E31").Value = "  +1.69%  "
Set-TextValue $ws.Range("D32") "4.150"
$ws.Range("E32").Value = "  +0.10%  "
Set-TextValue $ws.Range("D33") "4.084"
$ws.Range("E33").Value = "  -0.23%  "
Set-TextValue $ws.Range("D34") "1.910"
$ws.Range("E34").Value = "  +1.27%  "
Set-TextValue $ws.Range("D35") "1.186"
$ws.Range("E35").Value = "  +1.31%  "
Set-TextValue $ws.Range("D36") "0.7250"
$ws.Range("E36").Value = "  -0.10%  "
Set-TextValue $ws.Range("D37") "2.599"
$ws.Range("E37").Value = "  -0.39%  "
Set-TextValue $ws.Range("D38") "2.803"
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("E39").Value = "  +0.84%  "
$ws.Range("E40").Value = "  +0.02%  "
Set-TextValue $ws.Range("D41") "0.9150"
$ws.Range("E41").Value = "  -1.16%  "
Set-TextValue $ws.Range("D42") "6.137"
$ws.Range("E42").Value = "  -0.31%  "
$ws.Range("D43").Value = "2.040.06"
$ws.Range("E43").Value = "  +0.54%  "
Set-TextValue $ws.Range("D44") "0.9999"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D45") "67.66"
$ws.Range("E45").Value = "  +2.29%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-TextValue $ws.Range("D46") "101.95"
$ws.Range("E46").Value = "  -0.07%  "
Set-TextValue $ws.Range("D47") "7.507"
$ws.Range("E47").Value = "  +12.23%  "
Set-TextValue $ws.Range("D48") "0.4071"
$ws.Range("E48").Value = "  +0.29%  "
Set-TextValue $ws.Range("D49") "9.179"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue $ws.Range("D50") "1.686"
$ws.Range("E50").Value = "  +2.28%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D51") "0.1132"
$ws.Range("E51").Value = "  +1.20%  "
